$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> column letter -> new value, as per the target diff.
$updates = @{
    2 = @{ I = 0.5586792620790276; J = 0.5586792620790275;
           M = 46.29121633333333;  N = 138.873649;
           O = 0.3133663986859022; P = 0.3133663986859022;
           Q = 5.713601388779778;  R = 51.422412499018;
           S = 0.1750713083782022; T = 0.1750713083782022 }
    3 = @{ I = 0.5586792620790276; J = 0.5586792620790275;
           M = 46.81622333333333;
           O = 0.3169204109998198; P = 0.3169204109998198;
           Q = 5.778401602771111;  R = 52.00561442494;
           S = 0.1770568613551615; T = 0.1770568613551614 }
    4 = @{ I = 0.5586792620790276; J = 0.5586792620790275;
           M = 38.53544233333333;  N = 115.606327;
           O = 0.2608640200510233; P = 0.2608640200510233;
           Q = 4.756326886023778;  R = 42.806941974214;
           S = 0.1457393182250744; T = 0.1457393182250744 }
    5 = @{ I = 0.5586792620790276; J = 0.5586792620790275;
           M = 16.07945366666667;  N = 48.238361;
           O = 0.1088491702632547; P = 0.1088491702632547;
           Q = 1.984644087533555;  R = 17.861796787802;
           S = 0.06081177412058959; T = 0.06081177412058958 }
    6 = @{ G = 0.09749966666666667; H = 0.292499;
           I = 0.4413207379209724;  J = 0.4413207379209724;
           M = 46.29121633333333;   N = 138.873649;
           O = 0.3133663986859022;  P = 0.3133663986859022;
           Q = 4.513378162094555;   R = 40.620403458851;
           S = 0.1382950903077;     T = 0.1382950903077 }
    7 = @{ G = 0.09749966666666667; H = 0.292499;
           I = 0.4413207379209724;  J = 0.4413207379209724;
           M = 46.81622333333333;
           O = 0.3169204109998198;  P = 0.3169204109998198;
           Q = 4.564566169592222;   R = 41.08109552633;
           S = 0.1398635496446583;  T = 0.1398635496446583 }
    8 = @{ G = 0.09749966666666667; H = 0.292499;
           I = 0.4413207379209724;  J = 0.4413207379209724;
           M = 38.53544233333333;   N = 115.606327;
           O = 0.2608640200510233;  P = 0.2608640200510233;
           Q = 3.757192782352555;   R = 33.814735041173;
           S = 0.115124701825949;   T = 0.1151247018259489 }
    9 = @{ G = 0.09749966666666667; H = 0.292499;
           I = 0.4413207379209724;  J = 0.4413207379209724;
           M = 16.07945366666667;   N = 48.238361;
           O = 0.1088491702632547;  P = 0.1088491702632547;
           Q = 1.567741372682111;   R = 14.109672354139;
           S = 0.04803739614266515; T = 0.04803739614266514 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
